$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1621.9259
$ws.Range("I40").Value = 1415.8334
$ws.Range("J40").Value = 1786.8
$ws.Range("K40").Value = 1415.8334
$ws.Range("L40").Value = 1786.8
$ws.Range("M40").Value = -1240.8334
$ws.Range("N40").Value = -2136.8
$ws.Range("H43").Value = 1529.4
$ws.Range("I43").Value = 1468.2
$ws.Range("J43").Value = 1560
$ws.Range("K43").Value = 1468.2
$ws.Range("L43").Value = 1560
$ws.Range("M43").Value = -1399.2
$ws.Range("N43").Value = -1698
$ws.Range("H62").Value = 1112.4736
$ws.Range("I62").Value = 1012.7273
$ws.Range("J62").Value = 1249.625
$ws.Range("K62").Value = 1012.7273
$ws.Range("L62").Value = 1249.625
$ws.Range("M62").Value = -388.7273
$ws.Range("N62").Value = -2497.625
$ws.Range("H64").Value = 3488.92
$ws.Range("I64").Value = 3086.6667
$ws.Range("J64").Value = 3543.7727
$ws.Range("K64").Value = 3086.6667
$ws.Range("L64").Value = 3543.7727
$ws.Range("M64").Value = -2838.6667
$ws.Range("N64").Value = -4039.7727
$ws.Range("H65").Value = 1112.4736
$ws.Range("I65").Value = 1012.7273
$ws.Range("J65").Value = 1249.625
$ws.Range("K65").Value = 5063.636500000001
$ws.Range("L65").Value = 6248.125
$ws.Range("M65").Value = -1943.636500000001
$ws.Range("N65").Value = -12488.125
$ws.Range("H67").Value = 3488.92
$ws.Range("I67").Value = 3086.6667
$ws.Range("J67").Value = 3543.7727
$ws.Range("K67").Value = 3086.6667
$ws.Range("L67").Value = 3543.7727
$ws.Range("M67").Value = -2228.6667
$ws.Range("N67").Value = -5259.7727
$ws.Range("H113").Value = 197493.06
$ws.Range("I113").Value = 287184.9
$ws.Range("J113").Value = 3160.7778
$ws.Range("K113").Value = 287184.9
$ws.Range("L113").Value = 3160.7778
$ws.Range("M113").Value = -283930.9
$ws.Range("N113").Value = -9668.7778
$ws.Range("H116").Value = 2386585.8
$ws.Range("I116").Value = 10205995
$ws.Range("J116").Value = 6765.478
$ws.Range("K116").Value = 10205995
$ws.Range("L116").Value = 6765.478
$ws.Range("M116").Value = -10202553
$ws.Range("N116").Value = -13649.478
$ws.Range("H125").Value = 6294.731
$ws.Range("I125").Value = 7201.222
$ws.Range("J125").Value = 5814.8237
$ws.Range("K125").Value = 64810.998
$ws.Range("L125").Value = 52333.4133
$ws.Range("M125").Value = -62350.998
$ws.Range("N125").Value = -57253.4133
$ws.Range("H129").Value = 1910.973
$ws.Range("I129").Value = 453.2857
$ws.Range("J129").Value = 2251.1
$ws.Range("K129").Value = 1359.8571
$ws.Range("L129").Value = 6753.299999999999
$ws.Range("M129").Value = 3640.1429
$ws.Range("N129").Value = -16753.3
$ws.Range("H132").Value = 2059219.2
$ws.Range("I132").Value = 2924838.2
$ws.Range("J132").Value = 3374.5
$ws.Range("K132").Value = 8774514.600000001
$ws.Range("L132").Value = 10123.5
$ws.Range("M132").Value = -8771984.600000001
$ws.Range("N132").Value = -15183.5
$ws.Range("H137").Value = 1117.2858
$ws.Range("I137").Value = 763.3333
$ws.Range("J137").Value = 2002.1666
$ws.Range("K137").Value = 2289.9999
$ws.Range("L137").Value = 6006.4998
$ws.Range("M137").Value = 260.0001000000002
$ws.Range("N137").Value = -11106.4998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1268.1052
$ws.Range("I61").Value = 728.1429000000001
$ws.Range("J61").Value = 2780
$ws.Range("K61").Value = 728.1429000000001
$ws.Range("L61").Value = 2780
$ws.Range("M61").Value = -516.1429000000001
$ws.Range("N61").Value = -3204
$ws.Range("H63").Value = 2302.0483
$ws.Range("I63").Value = 2298.7834
$ws.Range("J63").Value = 2400
$ws.Range("K63").Value = 2298.7834
$ws.Range("L63").Value = 2400
$ws.Range("M63").Value = -1612.7834
$ws.Range("N63").Value = -3772
$ws.Range("H66").Value = 2302.0483
$ws.Range("I66").Value = 2298.7834
$ws.Range("J66").Value = 2400
$ws.Range("K66").Value = 11493.917
$ws.Range("L66").Value = 12000
$ws.Range("M66").Value = -8061.916999999999
$ws.Range("N66").Value = -18864
$ws.Range("H102").Value = 1569.3103
$ws.Range("I102").Value = 1571.7858
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 1571.7858
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = 50.21419999999989
$ws.Range("N102").Value = -4744
$ws.Range("H132").Value = 2081.5312
$ws.Range("I132").Value = 1430.4546
$ws.Range("J132").Value = 3513.9
$ws.Range("K132").Value = 4291.3638
$ws.Range("L132").Value = 10541.7
$ws.Range("M132").Value = -1761.3638
$ws.Range("N132").Value = -15601.7
$ws.Range("H136").Value = 1268.1052
$ws.Range("I136").Value = 728.1429000000001
$ws.Range("J136").Value = 2780
$ws.Range("K136").Value = 2184.4287
$ws.Range("L136").Value = 8340
$ws.Range("M136").Value = 365.5712999999996
$ws.Range("N136").Value = -13440

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1827.6364
$ws.Range("I86").Value = 1502
$ws.Range("K86").Value = 1502
$ws.Range("M86").Value = -379
$ws.Range("H89").Value = 1827.6364
$ws.Range("I89").Value = 1502
$ws.Range("K89").Value = 7510
$ws.Range("M89").Value = -1894
$ws.Range("H105").Value = 948705.5
$ws.Range("I105").Value = 1422028.2
$ws.Range("J105").Value = 2060
$ws.Range("K105").Value = 1422028.2
$ws.Range("L105").Value = 2060
$ws.Range("M105").Value = -1420281.2
$ws.Range("N105").Value = -5554

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2400.9607
$ws.Range("I31").Value = 2085.6428
$ws.Range("J31").Value = 2784.8262
$ws.Range("K31").Value = 2085.6428
$ws.Range("L31").Value = 2784.8262
$ws.Range("M31").Value = -1790.6428
$ws.Range("N31").Value = -3374.8262
$ws.Range("H34").Value = 2400.9607
$ws.Range("I34").Value = 2085.6428
$ws.Range("J34").Value = 2784.8262
$ws.Range("K34").Value = 2085.6428
$ws.Range("L34").Value = 2784.8262
$ws.Range("M34").Value = -1883.6428
$ws.Range("N34").Value = -3188.8262

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 6042.5835
$ws.Range("I133").Value = 5485.1665
$ws.Range("J133").Value = 6600
$ws.Range("K133").Value = 16455.4995
$ws.Range("L133").Value = 19800
$ws.Range("M133").Value = -11395.4995
$ws.Range("N133").Value = -29920

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2436.875
$ws.Range("I80").Value = 2365
$ws.Range("J80").Value = 2480
$ws.Range("K80").Value = 2365
$ws.Range("L80").Value = 2480
$ws.Range("M80").Value = -1367
$ws.Range("N80").Value = -4476
$ws.Range("H83").Value = 2436.875
$ws.Range("I83").Value = 2365
$ws.Range("J83").Value = 2480
$ws.Range("K83").Value = 11825
$ws.Range("L83").Value = 12400
$ws.Range("M83").Value = -6833
$ws.Range("N83").Value = -22384
$ws.Range("H97").Value = 606.25
$ws.Range("I97").Value = 350
$ws.Range("K97").Value = 350
$ws.Range("M97").Value = 146
$ws.Range("H122").Value = 2054.5
$ws.Range("I122").Value = 1189.25
$ws.Range("J122").Value = 2631.3333
$ws.Range("K122").Value = 3567.75
$ws.Range("L122").Value = 7893.999899999999
$ws.Range("M122").Value = -1117.75
$ws.Range("N122").Value = -12793.9999
$ws.Range("H126").Value = 1298.2727
$ws.Range("I126").Value = 1161.5
$ws.Range("J126").Value = 1663
$ws.Range("K126").Value = 3484.5
$ws.Range("L126").Value = 4989
$ws.Range("M126").Value = -1014.5
$ws.Range("N126").Value = -9929

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2261.8333
$ws.Range("I40").Value = 2210.8235
$ws.Range("J40").Value = 2385.7144
$ws.Range("K40").Value = 2210.8235
$ws.Range("L40").Value = 2385.7144
$ws.Range("M40").Value = -2074.8235
$ws.Range("N40").Value = -2657.7144
$ws.Range("H46").Value = 1060
$ws.Range("I46").Value = 1200
$ws.Range("J46").Value = 966.6667
$ws.Range("K46").Value = 1200
$ws.Range("L46").Value = 966.6667
$ws.Range("M46").Value = -1012
$ws.Range("N46").Value = -1342.6667
$ws.Range("H68").Value = 8030.6665
$ws.Range("I68").Value = 10154.546
$ws.Range("J68").Value = 2190
$ws.Range("K68").Value = 10154.546
$ws.Range("L68").Value = 2190
$ws.Range("M68").Value = -9405.546
$ws.Range("N68").Value = -3688
$ws.Range("H71").Value = 8030.6665
$ws.Range("I71").Value = 10154.546
$ws.Range("J71").Value = 2190
$ws.Range("K71").Value = 50772.73
$ws.Range("L71").Value = 10950
$ws.Range("M71").Value = -47028.73
$ws.Range("N71").Value = -18438
$ws.Range("H82").Value = 1998.0358
$ws.Range("I82").Value = 1435.4667
$ws.Range("J82").Value = 2647.1538
$ws.Range("K82").Value = 1435.4667
$ws.Range("L82").Value = 2647.1538
$ws.Range("M82").Value = -1074.4667
$ws.Range("N82").Value = -3369.1538
$ws.Range("H85").Value = 1998.0358
$ws.Range("I85").Value = 1435.4667
$ws.Range("J85").Value = 2647.1538
$ws.Range("K85").Value = 1435.4667
$ws.Range("L85").Value = 2647.1538
$ws.Range("M85").Value = -187.4666999999999
$ws.Range("N85").Value = -5143.1538
$ws.Range("H100").Value = 2096.8823
$ws.Range("I100").Value = 1631
$ws.Range("J100").Value = 2511
$ws.Range("K100").Value = 1631
$ws.Range("L100").Value = 2511
$ws.Range("M100").Value = -1090
$ws.Range("N100").Value = -3593

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2759.2222
$ws.Range("I126").Value = 1881.6923
$ws.Range("J126").Value = 5040.8
$ws.Range("K126").Value = 5645.0769
$ws.Range("L126").Value = 15122.4
$ws.Range("M126").Value = -3175.0769
$ws.Range("N126").Value = -20062.4
$ws.Range("H132").Value = 1639.6207
$ws.Range("I132").Value = 730.7
$ws.Range("J132").Value = 2118
$ws.Range("K132").Value = 2192.1
$ws.Range("L132").Value = 6354
$ws.Range("M132").Value = 337.8999999999996
$ws.Range("N132").Value = -11414
$ws.Range("H136").Value = 1142.6296
$ws.Range("I136").Value = 875.3889
$ws.Range("J136").Value = 1677.1111
$ws.Range("K136").Value = 2626.1667
$ws.Range("L136").Value = 5031.3333
$ws.Range("M136").Value = -76.16670000000022
$ws.Range("N136").Value = -10131.3333
